$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.244.70"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "1.442.52"
$ws.Range("E3").Value = "  +2.73%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9184"
$ws.Range("E5").Value = "  -8.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.47"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3636"
$ws.Range("E7").Value = "  -1.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3078"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.75"
$ws.Range("E9").Value = "  -1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.020"
$ws.Range("E10").Value = "  +1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06489"
$ws.Range("E11").Value = "  +0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.326"
$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.42"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.031"
$ws.Range("E15").Value = "  -1.89%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001008"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "1.441.36"
$ws.Range("E17").Value = "  +2.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9353"
$ws.Range("E18").Value = "  -6.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05601"
$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.55"
$ws.Range("E20").Value = "  -3.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.359"
$ws.Range("E21").Value = "  -3.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.17"
$ws.Range("E22").Value = "  -3.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.76"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.241"
$ws.Range("E24").Value = "  -1.57%  "

$ws.Range("D25").Value = "20.268.20"
$ws.Range("E25").Value = "  +1.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.09"
$ws.Range("E26").Value = "  +2.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.061"
$ws.Range("E27").Value = "  -7.71%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").Value = "1.594.39"
$ws.Range("E29").Value = "  +2.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.04"
$ws.Range("E30").Value = "  +0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.996"
$ws.Range("E31").Value = "  -2.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.826"
$ws.Range("E32").Value = "  -8.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7866"
$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07645"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.458"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05786"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134"
$ws.Range("E37").Value = "  +4.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.640"
$ws.Range("E38").Value = "  -3.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01985"
$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.14"
$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1841"
$ws.Range("E41").Value = "  -3.12%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9279"
$ws.Range("E42").Value = "  -7.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.021"
$ws.Range("E43").Value = "  -16.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5184"
$ws.Range("E44").Value = "  -1.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.473"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.70"
$ws.Range("E46").Value = "  -4.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.69"
$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5087"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.728"
$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06362"
$ws.Range("E50").Value = "  +3.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9920"
$ws.Range("E51").Value = "  -1.04%  "
